# Update "想去人数" (interested-count) figures in both the "展览" sheet and
# the "全部类型" sheet to match the freshly scraped output.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3121
$ws1.Range("F7").Value = 267
$ws1.Range("F9").Value = 1105
$ws1.Range("F10").Value = 15453
$ws1.Range("F11").Value = 223
$ws1.Range("F14").Value = 6094
$ws1.Range("F15").Value = 619
$ws1.Range("F31").Value = 10959

# Sheet "全部类型" (All types) — same events, shifted one row down
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3121
$ws4.Range("F8").Value = 267
$ws4.Range("F10").Value = 1105
$ws4.Range("F11").Value = 15453
$ws4.Range("F12").Value = 223
$ws4.Range("F15").Value = 6094
$ws4.Range("F16").Value = 619
$ws4.Range("F33").Value = 10959
